$wb = $excel.ActiveWorkbook

# Sheet ALC Row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 74.333336
$ws.Range("I8").Value = 74.333336
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 223.000008
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -84.00000800000001
$ws.Range("N8").ClearContents()

# Sheet ALC Row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 991.63635
$ws.Range("I41").Value = 217.5
$ws.Range("K41").Value = 217.5
$ws.Range("M41").Value = 222.5

# Sheet ARM Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11335.16
$ws.Range("I32").Value = 10765.792
$ws.Range("K32").Value = 10765.792
$ws.Range("M32").Value = -10478.792

# Sheet ARM Row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 237.75
$ws.Range("I97").Value = 237.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 237.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 258.25
$ws.Range("N97").ClearContents()

# Sheet ARM Row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4084.2307
$ws.Range("I110").Value = 3235.111
$ws.Range("J110").Value = 5994.75
$ws.Range("K110").Value = 3235.111
$ws.Range("L110").Value = 5994.75
$ws.Range("M110").Value = -1190.111
$ws.Range("N110").Value = -10084.75

# Sheet ARM Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 817.6
$ws.Range("J122").Value = 200
$ws.Range("L122").Value = 600
$ws.Range("N122").Value = -5500

# Sheet BSM Row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Sheet BSM Row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# Sheet BSM Row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 947.9231
$ws.Range("I107").Value = 947.9231
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 947.9231
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 972.0769
$ws.Range("N107").ClearContents()

# Sheet BSM Row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2255.375
$ws.Range("I134").Value = 1008.6
$ws.Range("K134").Value = 3025.8
$ws.Range("M134").Value = -490.8000000000002

# Sheet BSM Row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 35000
$ws.Range("I140").Value = 35000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 35000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -29820
$ws.Range("N140").ClearContents()

# Sheet CRP Row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630

# Sheet CRP Row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184

# Sheet CRP Row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5085
$ws.Range("I94").Value = 6005.5
$ws.Range("J94").Value = 4471.3335
$ws.Range("K94").Value = 6005.5
$ws.Range("L94").Value = 4471.3335
$ws.Range("M94").Value = -5554.5
$ws.Range("N94").Value = -5373.3335

# Sheet CUL Row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42.53846
$ws.Range("I2").Value = 15.375
$ws.Range("K2").Value = 92.25
$ws.Range("M2").Value = 20.75

# Sheet CUL Row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 172.71428
$ws.Range("J6").Value = 101
$ws.Range("L6").Value = 303
$ws.Range("N6").Value = -529

# Sheet CUL Row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1773
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1773
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 5319
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -5665

# Sheet CUL Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2071.4546
$ws.Range("I122").Value = 1958.8334
$ws.Range("J122").Value = 2206.6
$ws.Range("K122").Value = 17629.5006
$ws.Range("L122").Value = 19859.4
$ws.Range("M122").Value = -15179.5006
$ws.Range("N122").Value = -24759.4

# Sheet GSM Row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 27499.5
$ws.Range("J57").Value = 27499.5
$ws.Range("L57").Value = 27499.5
$ws.Range("N57").Value = -29139.5

# Sheet GSM Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9374.666999999999
$ws.Range("I122").Value = 11562.25
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 34686.75
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -32236.75
$ws.Range("N122").Value = -19898.5

# Sheet GSM Row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15435336
$ws.Range("I126").Value = 41153892
$ws.Range("K126").Value = 123461676
$ws.Range("M126").Value = -123459206

# Sheet LTW Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3054.6667
$ws.Range("I7").Value = 2945
$ws.Range("K7").Value = 2945
$ws.Range("M7").Value = -2833

# Sheet LTW Row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 494.5
$ws.Range("I9").Value = 526.3333
$ws.Range("K9").Value = 526.3333
$ws.Range("M9").Value = -302.3333

# Sheet LTW Row 14
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Sheet LTW Row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37499.75
$ws.Range("I40").Value = 24999.5
$ws.Range("K40").Value = 24999.5
$ws.Range("M40").Value = -24863.5

# Sheet LTW Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1370.5714
$ws.Range("I46").Value = 1197.6666
$ws.Range("J46").Value = 1500.25
$ws.Range("K46").Value = 1197.6666
$ws.Range("L46").Value = 1500.25
$ws.Range("M46").Value = -1009.6666
$ws.Range("N46").Value = -1876.25

# Sheet LTW Row 58
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 20000
$ws.Range("K58").Value = 20000
$ws.Range("M58").Value = -19740

# Sheet LTW Row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Sheet LTW Row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Sheet LTW Row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4025.913
$ws.Range("I100").Value = 2974.7
$ws.Range("J100").Value = 4834.5386
$ws.Range("K100").Value = 2974.7
$ws.Range("L100").Value = 4834.5386
$ws.Range("M100").Value = -2433.7
$ws.Range("N100").Value = -5916.5386

# Sheet LTW Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3054.6667
$ws.Range("I126").Value = 2945
$ws.Range("K126").Value = 8835
$ws.Range("M126").Value = -6365

# Sheet WVR Row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1376
$ws.Range("I21").Value = 507.8
$ws.Range("K21").Value = 507.8
$ws.Range("M21").Value = -272.8

# Sheet WVR Row 33
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 45000
$ws.Range("J33").Value = 45000
$ws.Range("L33").Value = 45000
$ws.Range("N33").Value = -45500

# Sheet WVR Row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 1376
$ws.Range("I35").Value = 507.8
$ws.Range("K35").Value = 507.8
$ws.Range("M35").Value = -217.8

# Sheet WVR Row 36
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 45000
$ws.Range("J36").Value = 45000
$ws.Range("L36").Value = 45000
$ws.Range("N36").Value = -45500

# Sheet WVR Row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 18897.5
$ws.Range("I70").Value = 12795
$ws.Range("K70").Value = 12795
$ws.Range("M70").Value = -12480

# Sheet WVR Row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 18897.5
$ws.Range("I73").Value = 12795
$ws.Range("K73").Value = 12795
$ws.Range("M73").Value = -11703

# Sheet WVR Row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1143.1818
$ws.Range("I113").Value = 825
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 2475
$ws.Range("L113").Value = 5100
$ws.Range("M113").Value = -305
$ws.Range("N113").Value = -9440

# Sheet WVR Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 40856.54
$ws.Range("I136").Value = 46921.363
$ws.Range("K136").Value = 140764.089
$ws.Range("M136").Value = -138214.089
